# TCMuscle.xlsx edit: add a "Y" mark in column E for the existing sample
# rows (rows 2-64), and append new Maurolicus muelleri rows (77-93) with
# their trawl counts, species name, sample codes and a "Y" in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill column E ("Run?") with "Y" for rows 2 through 64 -------------
$ws.Range("E2:E64").Value = "Y"

# --- Append new rows for Maurolicus muelleri samples -------------------
$speciesName = "Maurolicus muelleri"

$newRow = 77

# Rows 77-86: trawl count 8, codes TCMM001-TCMM010
for ($n = 1; $n -le 10; $n++) {
    $code = "TCMM{0:D3}" -f $n

    $ws.Cells.Item($newRow, 1).Value = 8
    $ws.Cells.Item($newRow, 2).Value = $speciesName
    $ws.Cells.Item($newRow, 2).Font.Italic = $true
    $ws.Cells.Item($newRow, 3).Value = $code
    $ws.Cells.Item($newRow, 4).Value = "Y"

    $newRow++
}

# Rows 87-93: trawl count 4, codes TCMM157-TCMM163
for ($n = 157; $n -le 163; $n++) {
    $code = "TCMM{0:D3}" -f $n

    $ws.Cells.Item($newRow, 1).Value = 4
    $ws.Cells.Item($newRow, 2).Value = $speciesName
    $ws.Cells.Item($newRow, 2).Font.Italic = $true
    $ws.Cells.Item($newRow, 3).Value = $code
    $ws.Cells.Item($newRow, 4).Value = "Y"

    $newRow++
}

# --- Restore the view state (scrolled / selected range) ----------------
[void]$ws.Range("B86:B93").Select()
